$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8 (shifts existing rows 8-17 down to 9-18)
$ws.Rows(8).Insert()

# Make sure the newly inserted row has no inherited formatting (plain style)
$ws.Range("A8:B8").Style = "Normal"

# Populate the new "Z1" tag/feature_name row (init_values left blank)
$ws.Range("A8").Value = "Z1"
$ws.Range("B8").Value = "Z1"

# Update the active selection to match the saved view state
$ws.Range("F13").Select() | Out-Null
